# Update "want to go" counts (F column) on the "展览" (exhibition) sheet
# and the "全部类型" (all types) sheet to reflect newly generated data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet): F2, F3, F4
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 979
$wsExhibition.Range("F3").Value = 1977
$wsExhibition.Range("F4").Value = 440

# Sheet "全部类型" (4th sheet): F4, F5, F6
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 979
$wsAll.Range("F5").Value = 1977
$wsAll.Range("F6").Value = 440
